# Auto-generated edit script applying numeric updates to the Maduin_Profits workbook
# per the target diff. Values are set as literals (matching the committed OOXML).

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 927
$ws.Range("J17").Value = 963
$ws.Range("L17").Value = 2889
$ws.Range("N17").Value = -3225
$ws.Range("H33").Value = 284.6
$ws.Range("I33").Value = 139.57143
$ws.Range("K33").Value = 139.57143
$ws.Range("M33").Value = 89.42857000000001
$ws.Range("H40").Value = 2139.9
$ws.Range("I40").Value = 2100
$ws.Range("K40").Value = 2100
$ws.Range("M40").Value = -1925
$ws.Range("H41").Value = 328.25
$ws.Range("I41").Value = 303.7143
$ws.Range("K41").Value = 303.7143
$ws.Range("M41").Value = 136.2857
$ws.Range("H116").Value = 3200
$ws.Range("I116").Value = 3200
$ws.Range("K116").Value = 3200
$ws.Range("M116").Value = 242
$ws.Range("H129").Value = 3982.3333
$ws.Range("I129").Value = 5473.5
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 16420.5
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = -11420.5
$ws.Range("N129").Value = -13000
$ws.Range("H132").Value = 2334
$ws.Range("I132").Value = 2000.8
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6002.4
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3472.4
$ws.Range("N132").Value = -17060
$ws.Range("H137").Value = 812.9
$ws.Range("I137").Value = 767.44446
$ws.Range("K137").Value = 2302.33338
$ws.Range("M137").Value = 247.66662
$ws.Range("H138").Value = 3139.4
$ws.Range("J138").Value = 3399.5293
$ws.Range("L138").Value = 10198.5879
$ws.Range("N138").Value = -20478.5879

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2152.5454
$ws.Range("I2").Value = 1368
$ws.Range("K2").Value = 1368
$ws.Range("M2").Value = -1255
$ws.Range("H32").Value = 2896.138
$ws.Range("I32").Value = 2650.5
$ws.Range("K32").Value = 2650.5
$ws.Range("M32").Value = -2363.5
$ws.Range("H116").Value = 2152.5454
$ws.Range("I116").Value = 1368
$ws.Range("K116").Value = 1368
$ws.Range("M116").Value = 926

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2152.5454
$ws.Range("I3").Value = 1368
$ws.Range("K3").Value = 1368
$ws.Range("M3").Value = -1254
$ws.Range("H105").Value = 3898.5715
$ws.Range("I105").Value = 3298.3333
$ws.Range("J105").Value = 7500
$ws.Range("K105").Value = 3298.3333
$ws.Range("L105").Value = 7500
$ws.Range("M105").Value = -1551.3333
$ws.Range("N105").Value = -10994
$ws.Range("H134").Value = 1403.6666
$ws.Range("I134").Value = 1403.6666
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4210.9998
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1675.9998
$ws.Range("N134").ClearContents()

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 438.625
$ws.Range("I16").Value = 438.625
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 438.625
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -151.625
$ws.Range("N16").ClearContents()
$ws.Range("H86").Value = 12332761
$ws.Range("I86").Value = 14539353
$ws.Range("K86").Value = 14539353
$ws.Range("M86").Value = -14538230
$ws.Range("H89").Value = 12332761
$ws.Range("I89").Value = 14539353
$ws.Range("K89").Value = 72696765
$ws.Range("M89").Value = -72691149
$ws.Range("H107").Value = 345.05884
$ws.Range("I107").Value = 308.76923
$ws.Range("K107").Value = 308.76923
$ws.Range("M107").Value = 1611.23077
$ws.Range("H113").Value = 438.625
$ws.Range("I113").Value = 438.625
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 438.625
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1731.375
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2730.3333
$ws.Range("J122").Value = 2400
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100
$ws.Range("H141").Value = 108542.29
$ws.Range("I141").Value = 51332
$ws.Range("J141").Value = 151450
$ws.Range("K141").Value = 51332
$ws.Range("L141").Value = 151450
$ws.Range("M141").Value = -46152
$ws.Range("N141").Value = -161810

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3499
$ws.Range("J51").Value = 3499
$ws.Range("L51").Value = 10497
$ws.Range("N51").Value = -11417
$ws.Range("H134").Value = 4750
$ws.Range("I134").Value = 4750
$ws.Range("K134").Value = 14250
$ws.Range("M134").Value = -9180

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H40").Value = 2937
$ws.Range("I40").Value = 2937
$ws.Range("K40").Value = 2937
$ws.Range("M40").Value = -2801
$ws.Range("H69").Value = 7163
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7163
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 7163
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -8785
$ws.Range("H72").Value = 7163
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7163
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 21489
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -29601
$ws.Range("H74").Value = 32499.75
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 39999.668
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 39999.668
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -41995.668
$ws.Range("H77").Value = 32499.75
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 39999.668
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 119999.004
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -129983.004
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2254.25
$ws.Range("I126").Value = 1183.1111
$ws.Range("K126").Value = 3549.3333
$ws.Range("M126").Value = -1079.3333
